$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.02406616871877758
$ws.Range("J2").Value = 0.02406616871877757
$ws.Range("M2").Value = 0.9949089999999999
$ws.Range("N2").Value = 2.984727
$ws.Range("O2").Value = 0.1476822527339178
$ws.Range("P2").Value = 0.1476822527339178
$ws.Range("Q2").Value = 0.04459049483466666
$ws.Range("R2").Value = 0.401314453512
$ws.Range("S2").Value = 0.003554146011063617
$ws.Range("T2").Value = 0.003554146011063616
$ws.Range("I3").Value = 0.02406616871877758
$ws.Range("J3").Value = 0.02406616871877757
$ws.Range("O3").Value = 0.2453919293791607
$ws.Range("P3").Value = 0.2453919293791607
$ws.Range("S3").Value = 0.005905643574665235
$ws.Range("T3").Value = 0.005905643574665234
$ws.Range("I4").Value = 0.02406616871877758
$ws.Range("J4").Value = 0.02406616871877757
$ws.Range("M4").Value = 1.748891
$ws.Range("N4").Value = 5.246673
$ws.Range("O4").Value = 0.2596017954064887
$ws.Range("P4").Value = 0.2596017954064887
$ws.Range("Q4").Value = 0.07838296276533334
$ws.Range("R4").Value = 0.705446664888
$ws.Range("S4").Value = 0.006247620607950135
$ws.Range("T4").Value = 0.006247620607950134
$ws.Range("I5").Value = 0.02406616871877758
$ws.Range("J5").Value = 0.02406616871877757
$ws.Range("M5").Value = 0.7268083333333334
$ws.Range("N5").Value = 2.180425
$ws.Range("O5").Value = 0.1078859392893731
$ws.Range("P5").Value = 0.1078859392893731
$ws.Range("Q5").Value = 0.03257458042222222
$ws.Range("R5").Value = 0.2931712238
$ws.Range("S5").Value = 0.002596401217321848
$ws.Range("T5").Value = 0.002596401217321848
$ws.Range("I6").Value = 0.02406616871877758
$ws.Range("J6").Value = 0.02406616871877757
$ws.Range("M6").Value = 1.613051666666667
$ws.Range("N6").Value = 4.839155
$ws.Range("O6").Value = 0.2394380831910597
$ws.Range("P6").Value = 0.2394380831910597
$ws.Range("Q6").Value = 0.07229482496444443
$ws.Range("R6").Value = 0.65065342468
$ws.Range("S6").Value = 0.005762357307776743
$ws.Range("T6").Value = 0.005762357307776743
$ws.Range("I7").Value = 0.8626970447097064
$ws.Range("J7").Value = 0.8626970447097063
$ws.Range("M7").Value = 0.9949089999999999
$ws.Range("N7").Value = 2.984727
$ws.Range("O7").Value = 0.1476822527339178
$ws.Range("P7").Value = 0.1476822527339178
$ws.Range("Q7").Value = 1.598430085217333
$ws.Range("R7").Value = 14.385870766956
$ws.Range("S7").Value = 0.1274050429896229
$ws.Range("T7").Value = 0.1274050429896228
$ws.Range("I8").Value = 0.8626970447097064
$ws.Range("J8").Value = 0.8626970447097063
$ws.Range("O8").Value = 0.2453919293791607
$ws.Range("P8").Value = 0.2453919293791607
$ws.Range("S8").Value = 0.211698892271015
$ws.Range("T8").Value = 0.2116988922710149
$ws.Range("I9").Value = 0.8626970447097064
$ws.Range("J9").Value = 0.8626970447097063
$ws.Range("M9").Value = 1.748891
$ws.Range("N9").Value = 5.246673
$ws.Range("O9").Value = 0.2596017954064887
$ws.Range("P9").Value = 0.2596017954064887
$ws.Range("Q9").Value = 2.809784603582667
$ws.Range("R9").Value = 25.288061432244
$ws.Range("S9").Value = 0.2239577016985116
$ws.Range("T9").Value = 0.2239577016985116
$ws.Range("I10").Value = 0.8626970447097064
$ws.Range("J10").Value = 0.8626970447097063
$ws.Range("M10").Value = 0.7268083333333334
$ws.Range("N10").Value = 2.180425
$ws.Range("O10").Value = 0.1078859392893731
$ws.Range("P10").Value = 0.1078859392893731
$ws.Range("Q10").Value = 1.167697051877778
$ws.Range("R10").Value = 10.5092734669
$ws.Range("S10").Value = 0.09307288099067298
$ws.Range("T10").Value = 0.09307288099067297
$ws.Range("I11").Value = 0.8626970447097064
$ws.Range("J11").Value = 0.8626970447097063
$ws.Range("M11").Value = 1.613051666666667
$ws.Range("N11").Value = 4.839155
$ws.Range("O11").Value = 0.2394380831910597
$ws.Range("P11").Value = 0.2394380831910597
$ws.Range("Q11").Value = 2.591543862815555
$ws.Range("R11").Value = 23.32389476534
$ws.Range("S11").Value = 0.206562526759884
$ws.Range("T11").Value = 0.206562526759884
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.210882
$ws.Range("H12").Value = 0.6326459999999999
$ws.Range("I12").Value = 0.113236786571516
$ws.Range("J12").Value = 0.113236786571516
$ws.Range("M12").Value = 0.9949089999999999
$ws.Range("N12").Value = 2.984727
$ws.Range("O12").Value = 0.1476822527339178
$ws.Range("P12").Value = 0.1476822527339178
$ws.Range("Q12").Value = 0.209808399738
$ws.Range("R12").Value = 1.888275597642
$ws.Range("S12").Value = 0.01672306373323134
$ws.Range("T12").Value = 0.01672306373323134
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.210882
$ws.Range("H13").Value = 0.6326459999999999
$ws.Range("I13").Value = 0.113236786571516
$ws.Range("J13").Value = 0.113236786571516
$ws.Range("O13").Value = 0.2453919293791607
$ws.Range("P13").Value = 0.2453919293791607
$ws.Range("Q13").Value = 0.34862203859
$ws.Range("R13").Value = 3.13759834731
$ws.Range("S13").Value = 0.02778739353348056
$ws.Range("T13").Value = 0.02778739353348055
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.210882
$ws.Range("H14").Value = 0.6326459999999999
$ws.Range("I14").Value = 0.113236786571516
$ws.Range("J14").Value = 0.113236786571516
$ws.Range("M14").Value = 1.748891
$ws.Range("N14").Value = 5.246673
$ws.Range("O14").Value = 0.2596017954064887
$ws.Range("P14").Value = 0.2596017954064887
$ws.Range("Q14").Value = 0.368809631862
$ws.Range("R14").Value = 3.319286686758
$ws.Range("S14").Value = 0.02939647310002693
$ws.Range("T14").Value = 0.02939647310002693
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.210882
$ws.Range("H15").Value = 0.6326459999999999
$ws.Range("I15").Value = 0.113236786571516
$ws.Range("J15").Value = 0.113236786571516
$ws.Range("M15").Value = 0.7268083333333334
$ws.Range("N15").Value = 2.180425
$ws.Range("O15").Value = 0.1078859392893731
$ws.Range("P15").Value = 0.1078859392893731
$ws.Range("Q15").Value = 0.15327079495
$ws.Range("R15").Value = 1.37943715455
$ws.Range("S15").Value = 0.01221665708137828
$ws.Range("T15").Value = 0.01221665708137828
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.210882
$ws.Range("H16").Value = 0.6326459999999999
$ws.Range("I16").Value = 0.113236786571516
$ws.Range("J16").Value = 0.113236786571516
$ws.Range("M16").Value = 1.613051666666667
$ws.Range("N16").Value = 4.839155
$ws.Range("O16").Value = 0.2394380831910597
$ws.Range("P16").Value = 0.2394380831910597
$ws.Range("Q16").Value = 0.3401635615699999
$ws.Range("R16").Value = 3.06147205413
$ws.Range("S16").Value = 0.02711319912339892
$ws.Range("T16").Value = 0.02711319912339892

Write-Host "Applied 160 cell updates"